# Refresh NATMI LR-pair edge/expression metrics (Psen1-Notch1) with
# recomputed TPM-based values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 9.409481333333334
$ws.Range("H2").Value = 28.228444
$ws.Range("I2").Value = 0.2433300530093958
$ws.Range("J2").Value = 0.2433300530093958
$ws.Range("M2").Value = 48.42420966666666
$ws.Range("N2").Value = 145.272629
$ws.Range("O2").Value = 0.6311762527593259
$ws.Range("P2").Value = 0.6311762527593258
$ws.Range("Q2").Value = 455.6466969399195
$ws.Range("R2").Value = 4100.820272459277
$ws.Range("S2").Value = 0.1535841510421986
$ws.Range("T2").Value = 0.1535841510421986
$ws.Range("G3").Value = 9.409481333333334
$ws.Range("H3").Value = 28.228444
$ws.Range("I3").Value = 0.2433300530093958
$ws.Range("J3").Value = 0.2433300530093958
$ws.Range("M3").Value = 6.849914666666667
$ws.Range("O3").Value = 0.08928392431779728
$ws.Range("P3").Value = 0.08928392431779726
$ws.Range("Q3").Value = 64.45414419092623
$ws.Range("R3").Value = 580.0872977183361
$ws.Range("S3").Value = 0.02172546203713649
$ws.Range("T3").Value = 0.02172546203713649
$ws.Range("G4").Value = 9.409481333333334
$ws.Range("H4").Value = 28.228444
$ws.Range("I4").Value = 0.2433300530093958
$ws.Range("J4").Value = 0.2433300530093958
$ws.Range("N4").Value = 64.33937399999999
$ws.Range("O4").Value = 0.2795398229228769
$ws.Range("P4").Value = 0.2795398229228769
$ws.Range("Q4").Value = 201.8000462171173
$ws.Range("R4").Value = 1816.200415954056
$ws.Range("S4").Value = 0.06802043993006075
$ws.Range("T4").Value = 0.06802043993006075
$ws.Range("I5").Value = 0.5069354697952918
$ws.Range("J5").Value = 0.5069354697952919
$ws.Range("M5").Value = 48.42420966666666
$ws.Range("N5").Value = 145.272629
$ws.Range("O5").Value = 0.6311762527593259
$ws.Range("P5").Value = 0.6311762527593258
$ws.Range("Q5").Value = 949.2599435096988
$ws.Range("R5").Value = 8543.33949158729
$ws.Range("S5").Value = 0.3199656302161807
$ws.Range("T5").Value = 0.3199656302161807
$ws.Range("I6").Value = 0.5069354697952918
$ws.Range("J6").Value = 0.5069354697952919
$ws.Range("M6").Value = 6.849914666666667
$ws.Range("O6").Value = 0.08928392431779728
$ws.Range("P6").Value = 0.08928392431779726
$ws.Range("Q6").Value = 134.2789000437156
$ws.Range("S6").Value = 0.04526118811920984
$ws.Range("T6").Value = 0.04526118811920984
$ws.Range("I7").Value = 0.5069354697952918
$ws.Range("J7").Value = 0.5069354697952919
$ws.Range("N7").Value = 64.33937399999999
$ws.Range("O7").Value = 0.2795398229228769
$ws.Range("P7").Value = 0.2795398229228769
$ws.Range("Q7").Value = 420.4149876621933
$ws.Range("R7").Value = 3783.73488895974
$ws.Range("S7").Value = 0.1417086514599013
$ws.Range("T7").Value = 0.1417086514599013
$ws.Range("G8").Value = 9.657138
$ws.Range("I8").Value = 0.2497344771953123
$ws.Range("J8").Value = 0.2497344771953124
$ws.Range("M8").Value = 48.42420966666666
$ws.Range("N8").Value = 145.272629
$ws.Range("O8").Value = 0.6311762527593259
$ws.Range("P8").Value = 0.6311762527593258
$ws.Range("Q8").Value = 467.6392752919339
$ws.Range("R8").Value = 4208.753477627406
$ws.Range("S8").Value = 0.1576264715009466
$ws.Range("T8").Value = 0.1576264715009466
$ws.Range("G9").Value = 9.657138
$ws.Range("I9").Value = 0.2497344771953123
$ws.Range("J9").Value = 0.2497344771953124
$ws.Range("M9").Value = 6.849914666666667
$ws.Range("O9").Value = 0.08928392431779728
$ws.Range("P9").Value = 0.08928392431779726
$ws.Range("Q9").Value = 66.150571224224
$ws.Range("S9").Value = 0.02229727416145093
$ws.Range("T9").Value = 0.02229727416145093
$ws.Range("G10").Value = 9.657138
$ws.Range("I10").Value = 0.2497344771953123
$ws.Range("J10").Value = 0.2497344771953124
$ws.Range("N10").Value = 64.33937399999999
$ws.Range("O10").Value = 0.2795398229228769
$ws.Range("P10").Value = 0.2795398229228769
$ws.Range("Q10").Value = 207.111404517204
$ws.Range("R10").Value = 1864.002640654836
$ws.Range("S10").Value = 0.06981073153291484
$ws.Range("T10").Value = 0.06981073153291485
